$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting rows 138:171 down to 139:172.
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new data record.
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C138").Value = 'Los Lagos'
$ws.Range("D138").Value = 44508
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 100112043
$ws.Range("G138").Value = 'Pepino ensalada'
$ws.Range("H138").Value = 'Sin especificar'
$ws.Range("I138").Value = 'Primera'
$ws.Range("J138").Value = 200
$ws.Range("K138").Value = 12000
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 12000
$ws.Range("N138").Value = '$/caja 60 unidades'
$ws.Range("O138").Value = 'Región de Arica y Parinacota'
$ws.Range("P138").Value = 200
$ws.Range("Q138").Value = 60
$ws.Range("R138").Value = 'Hortaliza'
